$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output ("ShapeCount=" + $ws.Shapes.Count)
$shp = $ws.Shapes.AddShape(1, 10, 10, 10, 10)
Write-Output ("ShapeCountAfterAdd=" + $ws.Shapes.Count)
$shp.Delete()
Write-Output ("ShapeCountAfterDelete=" + $ws.Shapes.Count)
$ws.Range("D1").Value = "res"
